$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C: decimal representation of the hex ROM bytes already present in
# column B ("We can check which of all sensors is online!").
$data = @(
    @(2,  "40 255 144 37 164 22 4 65"),
    @(3,  "40 255 32 150 164 22 4 23"),
    @(4,  "40 255 66 71 164 22 4 11"),
    @(5,  "40 255 161 155 164 22 5 210"),
    @(6,  "40 255 187 13 164 22 5 251"),
    @(7,  "40 255 118 161 164 22 5 142"),
    @(8,  "40 255 151 72 164 22 4 29"),
    @(9,  "40 255 21 11 164 22 5 99"),
    @(10, "40 255 91 156 164 22 5 79"),
    @(11, "40 255 114 96 164 22 4 121"),
    @(12, "40 255 74 86 164 22 4 130"),
    @(13, "40 255 17 78 164 22 4 67"),
    @(14, "40 255 43 131 164 22 4 222"),
    @(15, "40 255 7 148 164 22 4 185"),
    @(16, "40 255 178 12 164 22 5 135"),
    @(17, "40 255 118 224 148 22 4 97")
)

foreach ($item in $data) {
    $row = $item[0]
    $val = $item[1]
    $ws.Cells.Item($row, 3).Value = $val
}

# First data row got a "0" number format applied (kept as text) while the
# rest of the column stayed on the default/general style.
$ws.Cells.Item(2, 3).NumberFormat = "0"

$ws.Columns.Item(1).ColumnWidth = 17.5
$ws.Columns.Item(2).ColumnWidth = 31
$ws.Columns.Item(3).ColumnWidth = 31.5

$null = $ws.Range("C31").Select()

Write-Host "done"
